$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0006075818656279264
$ws.Range("C2").Value = 0.04103571897497393
$ws.Range("D2").Value = 189.6080260415259
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("G2").Value = 190.1830553009682

$ws.Range("B3").Value = 3.272327238179451
$ws.Range("C3").Value = 1.626987699542094
$ws.Range("D3").Value = 18.71679738969934
$ws.Range("E3").Value = 0.5333859586016987
$ws.Range("G3").Value = 24.14949828602258
